# Updates the cryptos price/volume table (columns D and E) with
# refreshed values from the latest scrape, matching the commit
# "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.325.78"
$ws.Range("E2").Value = "  -0.93%  "

$ws.Range("D3").Value = "3.238.46"
$ws.Range("E3").Value = "  +2.97%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.45%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "3.236.46"
$ws.Range("E8").Value = "  +3.08%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.53%  "

$ws.Range("E10").Value = "  -1.18%  "

$ws.Range("E11").Value = "  -0.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.464"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("E13").Value = "  -2.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.33"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").Value = "3.773.01"
$ws.Range("E15").Value = "  +2.88%  "

$ws.Range("E16").Value = "  -0.27%  "

$ws.Range("D17").Value = "3.240.64"
$ws.Range("E17").Value = "  +3.06%  "

$ws.Range("D18").Value = "63.342.56"
$ws.Range("E18").Value = "  -1.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.722"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.64%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.91"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.90%  "

$ws.Range("E28").Value = "  -1.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.05"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.01%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.31%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.88%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").Value = "  -2.71%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.52"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.89%  "

$ws.Range("E35").Value = "  -1.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.63%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.07%  "

$ws.Range("E38").Value = "  -3.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0392"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "421.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.69%  "

$ws.Range("D42").Value = "2.975.51"
$ws.Range("E42").Value = "  +1.38%  "

$ws.Range("E43").Value = "  -7.00%  "

$ws.Range("E44").Value = "  -8.06%  "

$ws.Range("E45").Value = "  +2.71%  "

$ws.Range("E46").Value = "  -1.88%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("E48").Value = "  -2.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("E50").Value = "  +0.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.25"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.81%  "
